$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "sequence/run_0629_0618_pooled/"

for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val.StartsWith($prefix)) {
        $cell.Value2 = $val.Substring($prefix.Length)
    }
}

$ws.Range("F2:F25").Select()
